# Using Page Object Model for efficiency
# Add the new "WishList_Addition" sheet right after "Password_Validation".
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "WishList_Addition"

# Pre-format the working range as Text (reuses the same style already present
# in the workbook - avoids creating new cellXfs entries).
$new.Range("A1:J9").NumberFormat = "@"

# Header row.
$new.Range("A1").Value = '${username}'
$new.Range("B1").Value = '${password}'
$new.Range("C1").Value = '${Items}'
$new.Range("D1").Value = '${Scroll}'

# Row data: username / password / item xpath / scroll command.
# Filled column-by-column (A, then B, then C, then D) to match the shared
# string intern order of the authored workbook.
$items = @(
    "xpath=//img[@alt='Duplex study table']",
    "xpath=//img[@alt='Dash wall mounted study table']",
    "xpath=//img[@alt='Tuck fold away work desk']",
    "xpath=//img[@alt='POD 180 large study table']",
    "xpath=//img[@alt='POD 180 small study table']",
    'xpath=//img[@alt="Fluid portable study table"]',
    'xpath=//img[@alt="Oblique study table"]',
    'xpath=//img[@alt="Step-up compact study table"]'
)

for ($row = 2; $row -le 9; $row++) {
    $new.Cells.Item($row, 1).Value = "cse20013@tezu.ac.in"
}
for ($row = 2; $row -le 9; $row++) {
    $new.Cells.Item($row, 2).Value = "usy253qu"
}
for ($i = 0; $i -lt $items.Length; $i++) {
    $new.Cells.Item($i + 2, 3).Value = $items[$i]
}
for ($row = 2; $row -le 9; $row++) {
    $new.Cells.Item($row, 4).Value = "window.scrollTo(0,200)"
}

# Hyperlink column A (rows 2-9) -- matches the style already used for
# mailto hyperlinks on the Password_Validation sheet.
for ($row = 2; $row -le 9; $row++) {
    $new.Hyperlinks.Add($new.Cells.Item($row, 1), "mailto:cse20013@tezu.ac.in") | Out-Null
}

# Hyperlinks.Add() re-styles the target cell with a fresh (unused) style;
# restore the correct Hyperlink+Text style by pasting the format already
# used by the existing hyperlink cells on sheet 1.
$ws1.Range("A2").Copy()
$new.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection on the new sheet matches the authored file.
$new.Range("E6").Select()
